$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Property1" to "DataNode"
$ws.Name = "DataNode"

# Move the active selection to D42 (matches the authored selection state)
$ws.Range("D42").Select()
